$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1186.25
$ws.Cells.Item(18, 9).Value = 927.1429000000001
$ws.Cells.Item(18, 10).Value = 3000
$ws.Cells.Item(18, 11).Value = 927.1429000000001
$ws.Cells.Item(18, 12).Value = 3000
$ws.Cells.Item(18, 13).Value = -643.1429000000001
$ws.Cells.Item(18, 14).Value = -3568
$ws.Cells.Item(41, 8).Value = 411.75
$ws.Cells.Item(41, 9).Value = 391.92856
$ws.Cells.Item(41, 11).Value = 391.92856
$ws.Cells.Item(41, 13).Value = 48.07144
$ws.Cells.Item(62, 8).Value = 3226.4
$ws.Cells.Item(62, 9).Value = 3083.3333
$ws.Cells.Item(62, 10).Value = 3441
$ws.Cells.Item(62, 11).Value = 3083.3333
$ws.Cells.Item(62, 12).Value = 3441
$ws.Cells.Item(62, 13).Value = -2459.3333
$ws.Cells.Item(62, 14).Value = -4689
$ws.Cells.Item(65, 8).Value = 3226.4
$ws.Cells.Item(65, 9).Value = 3083.3333
$ws.Cells.Item(65, 10).Value = 3441
$ws.Cells.Item(65, 11).Value = 15416.6665
$ws.Cells.Item(65, 12).Value = 17205
$ws.Cells.Item(65, 13).Value = -12296.6665
$ws.Cells.Item(65, 14).Value = -23445
$ws.Cells.Item(98, 8).Value = 583
$ws.Cells.Item(98, 9).Value = 583
$ws.Cells.Item(98, 11).Value = 583
$ws.Cells.Item(98, 13).Value = 915
$ws.Cells.Item(122, 8).Value = 583
$ws.Cells.Item(122, 9).Value = 583
$ws.Cells.Item(122, 11).Value = 1749
$ws.Cells.Item(122, 13).Value = 701
$ws.Cells.Item(126, 8).Value = 47000
$ws.Cells.Item(126, 10).Value = 47000
$ws.Cells.Item(126, 12).Value = 47000
$ws.Cells.Item(126, 14).Value = -56880
$ws.Cells.Item(135, 8).Value = 898.4545000000001
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 898.4545000000001
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 8086.0905
$ws.Cells.Item(135, 13).Value = ""
$ws.Cells.Item(135, 14).Value = -13156.0905
$ws.Cells.Item(140, 8).Value = 20000
$ws.Cells.Item(140, 9).Value = 20000
$ws.Cells.Item(140, 11).Value = 20000
$ws.Cells.Item(140, 13).Value = -14820

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 500.4
$ws.Cells.Item(12, 9).Value = 500.4
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 500.4
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = -327.4
$ws.Cells.Item(12, 14).Value = ""
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).Value = ""
$ws.Cells.Item(56, 14).Value = ""
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 13).Value = ""
$ws.Cells.Item(37, 8).Value = 1313
$ws.Cells.Item(37, 9).Value = 426
$ws.Cells.Item(37, 10).Value = 2200
$ws.Cells.Item(37, 11).Value = 426
$ws.Cells.Item(37, 12).Value = 2200
$ws.Cells.Item(37, 13).Value = -289
$ws.Cells.Item(37, 14).Value = -2474
$ws.Cells.Item(88, 8).Value = 32059.6
$ws.Cells.Item(88, 10).Value = 32059.6
$ws.Cells.Item(88, 12).Value = 32059.6
$ws.Cells.Item(88, 14).Value = -32871.6
$ws.Cells.Item(91, 8).Value = 32059.6
$ws.Cells.Item(91, 10).Value = 32059.6
$ws.Cells.Item(91, 12).Value = 32059.6
$ws.Cells.Item(91, 14).Value = -34867.6
$ws.Cells.Item(107, 8).Value = 135333
$ws.Cells.Item(107, 10).Value = 2999
$ws.Cells.Item(107, 12).Value = 2999
$ws.Cells.Item(107, 14).Value = -6839
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 5617
$ws.Cells.Item(5, 9).Value = 7705.25
$ws.Cells.Item(5, 10).Value = 2832.6667
$ws.Cells.Item(5, 11).Value = 7705.25
$ws.Cells.Item(5, 12).Value = 2832.6667
$ws.Cells.Item(5, 13).Value = -7593.25
$ws.Cells.Item(5, 14).Value = -3056.6667
$ws.Cells.Item(22, 8).Value = 645.6667
$ws.Cells.Item(22, 9).Value = 799.5714
$ws.Cells.Item(22, 10).Value = 430.2
$ws.Cells.Item(22, 11).Value = 799.5714
$ws.Cells.Item(22, 12).Value = 430.2
$ws.Cells.Item(22, 13).Value = -449.5714
$ws.Cells.Item(22, 14).Value = -1130.2
$ws.Cells.Item(25, 8).Value = 5000
$ws.Cells.Item(25, 10).Value = 5000
$ws.Cells.Item(25, 12).Value = 5000
$ws.Cells.Item(25, 14).Value = -5348
$ws.Cells.Item(53, 8).Value = 55536.8
$ws.Cells.Item(53, 10).Value = 55536.8
$ws.Cells.Item(53, 12).Value = 55536.8
$ws.Cells.Item(53, 14).Value = -56750.8
$ws.Cells.Item(59, 8).Value = 104
$ws.Cells.Item(59, 9).Value = 104
$ws.Cells.Item(59, 11).Value = 104
$ws.Cells.Item(59, 13).Value = 1041
$ws.Cells.Item(60, 8).Value = 20000
$ws.Cells.Item(60, 9).Value = 20000
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 20000
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = -19489
$ws.Cells.Item(60, 14).Value = ""
$ws.Cells.Item(96, 8).Value = 27532
$ws.Cells.Item(96, 10).Value = 27532
$ws.Cells.Item(96, 12).Value = 27532
$ws.Cells.Item(96, 14).Value = -33024
$ws.Cells.Item(99, 8).Value = 910535.4399999999
$ws.Cells.Item(99, 10).Value = 5000000
$ws.Cells.Item(99, 12).Value = 5000000
$ws.Cells.Item(99, 14).Value = -5002996
$ws.Cells.Item(107, 8).Value = 940
$ws.Cells.Item(107, 9).Value = 800
$ws.Cells.Item(107, 10).Value = 1500
$ws.Cells.Item(107, 11).Value = 800
$ws.Cells.Item(107, 12).Value = 1500
$ws.Cells.Item(107, 13).Value = 1120
$ws.Cells.Item(107, 14).Value = -5340
$ws.Cells.Item(126, 8).Value = 910535.4399999999
$ws.Cells.Item(126, 10).Value = 5000000
$ws.Cells.Item(126, 12).Value = 15000000
$ws.Cells.Item(126, 14).Value = -15004940

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 62116.668
$ws.Cells.Item(4, 9).Value = 1473.5714
$ws.Cells.Item(4, 11).Value = 4420.7142
$ws.Cells.Item(4, 13).Value = -4308.7142
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = ""
$ws.Cells.Item(9, 14).Value = ""
$ws.Cells.Item(92, 8).Value = 841
$ws.Cells.Item(92, 9).Value = 841
$ws.Cells.Item(92, 11).Value = 2523
$ws.Cells.Item(92, 13).Value = -1275
$ws.Cells.Item(104, 8).Value = 3010
$ws.Cells.Item(104, 10).Value = 4020
$ws.Cells.Item(104, 12).Value = 12060
$ws.Cells.Item(104, 14).Value = -17302
$ws.Cells.Item(128, 8).Value = 633332.7
$ws.Cells.Item(128, 9).Value = 633332.7
$ws.Cells.Item(128, 11).Value = 1899998.1
$ws.Cells.Item(128, 13).Value = -1895018.1
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 13).Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 13).Value = ""
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).Value = ""
$ws.Cells.Item(122, 8).Value = 5833.75
$ws.Cells.Item(122, 9).Value = 5699.7
$ws.Cells.Item(122, 11).Value = 17099.1
$ws.Cells.Item(122, 13).Value = -14649.1
$ws.Cells.Item(126, 8).Value = 1999.6666
$ws.Cells.Item(126, 9).Value = 1999.6666
$ws.Cells.Item(126, 11).Value = 5998.9998
$ws.Cells.Item(126, 13).Value = -3528.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 14333.333
$ws.Cells.Item(4, 9).Value = 5000
$ws.Cells.Item(4, 11).Value = 5000
$ws.Cells.Item(4, 13).Value = -4887
$ws.Cells.Item(22, 8).Value = 666.6667
$ws.Cells.Item(22, 9).Value = 500
$ws.Cells.Item(22, 10).Value = 750
$ws.Cells.Item(22, 11).Value = 500
$ws.Cells.Item(22, 12).Value = 750
$ws.Cells.Item(22, 13).Value = -205
$ws.Cells.Item(22, 14).Value = -1340
$ws.Cells.Item(27, 8).Value = 666.6667
$ws.Cells.Item(27, 9).Value = 500
$ws.Cells.Item(27, 10).Value = 750
$ws.Cells.Item(27, 11).Value = 500
$ws.Cells.Item(27, 12).Value = 750
$ws.Cells.Item(27, 13).Value = -393
$ws.Cells.Item(27, 14).Value = -964
$ws.Cells.Item(28, 8).Value = 14333.333
$ws.Cells.Item(28, 9).Value = 5000
$ws.Cells.Item(28, 11).Value = 5000
$ws.Cells.Item(28, 13).Value = -4768
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 14).Value = ""
$ws.Cells.Item(37, 8).Value = 14333.333
$ws.Cells.Item(37, 9).Value = 5000
$ws.Cells.Item(37, 11).Value = 5000
$ws.Cells.Item(37, 13).Value = -4893
$ws.Cells.Item(43, 8).Value = 10000
$ws.Cells.Item(43, 10).Value = 10000
$ws.Cells.Item(43, 12).Value = 10000
$ws.Cells.Item(43, 14).Value = -10386
$ws.Cells.Item(94, 8).Value = 44662
$ws.Cells.Item(94, 10).Value = 44662
$ws.Cells.Item(94, 12).Value = 44662
$ws.Cells.Item(94, 14).Value = -46014

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1347
$ws.Cells.Item(122, 9).Value = 1096.25
$ws.Cells.Item(122, 11).Value = 3288.75
$ws.Cells.Item(122, 13).Value = -838.75
$ws.Cells.Item(126, 8).Value = 4033.5454
$ws.Cells.Item(126, 9).Value = 3436.4
$ws.Cells.Item(126, 11).Value = 10309.2
$ws.Cells.Item(126, 13).Value = -7839.200000000001
